$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain plain text while we assign the new
# values below - many of the new quotes are single-dot decimals (e.g. "242.56")
# that Excel would otherwise auto-convert to numbers, destroying the original
# "29.858.75"-style formatting used throughout this sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.853.65"
$ws.Range("D3").Value = "1.886.76"
$ws.Range("D5").Value = "0.7450"
$ws.Range("D6").Value = "242.56"
$ws.Range("D8").Value = "0.3110"
$ws.Range("D9").Value = "25.33"
$ws.Range("D10").Value = "0.07108"
$ws.Range("D11").Value = "0.08486"
$ws.Range("D12").Value = "0.7596"
$ws.Range("D13").Value = "1.916.23"
$ws.Range("D14").Value = "5.348"
$ws.Range("D15").Value = "93.25"
$ws.Range("D16").Value = "6.144"
$ws.Range("D17").Value = "29.916.77"
$ws.Range("D18").Value = "13.69"
$ws.Range("D19").Value = "243.14"
$ws.Range("D20").Value = "0.000007788"
$ws.Range("D21").Value = "2.159.44"
$ws.Range("D22").Value = "1.000"
$ws.Range("D23").Value = "7.985"
$ws.Range("D24").Value = "1.002"
$ws.Range("D25").Value = "0.1584"
$ws.Range("D26").Value = "9.376"
$ws.Range("D27").Value = "162.38"
$ws.Range("D28").Value = "18.73"
$ws.Range("D30").Value = "1.506"
$ws.Range("D31").Value = "1.533"
$ws.Range("D32").Value = "4.465"
$ws.Range("D33").Value = "4.098"
$ws.Range("D34").Value = "0.05395"
$ws.Range("D35").Value = "1.238"
$ws.Range("D36").Value = "0.7430"
$ws.Range("D38").Value = "2.711"
$ws.Range("D39").Value = "0.01931"
$ws.Range("D40").Value = "2.768"
$ws.Range("D42").Value = "6.063"
$ws.Range("D44").Value = "1.085.97"
$ws.Range("D45").Value = "0.8626"
$ws.Range("D46").Value = "1.001"
$ws.Range("D47").Value = "102.57"
$ws.Range("D48").Value = "7.663"
$ws.Range("D49").Value = "1.860"
$ws.Range("D50").Value = "3.065"
$ws.Range("D51").Value = "2.055.33"

# Restore the default (unstyled) cell style so the text-format override above
# does not leave a stray style applied compared to the original workbook.
$ws.Range("D2:D51").Style = "Normal"

# Volume(1h) column (E) values are already unambiguous text (percent sign +
# surrounding spaces), so a plain Value assignment keeps them as strings.
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("E5").Value = "  -4.83%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("E11").Value = "  +4.69%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  -1.03%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("E44").Value = "  -5.04%  "
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("E51").Value = "  +0.96%  "

